# Apply updated odds values to Sheet1, as per the upstream FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Fortaleza - Atletico-MG)
$ws.Range("G2").Value = 1.73
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("X2").Value = 7
$ws.Range("AD2").Value = 7
$ws.Range("AN2").Value = 3.5
$ws.Range("AT2").Value = 2.5

# Row 5 (Sport Recife - Operario)
$ws.Range("G5").Value = 1.6
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.25
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("AD5").Value = 7.5
$ws.Range("AH5").Value = 12
$ws.Range("AO5").Value = 8.5
$ws.Range("AQ5").Value = 29

# Row 11 (Danubio - Penarol)
$ws.Range("AG11").Value = 1250

# Row 12 (Boston River - Progreso)
$ws.Range("I12").Value = 3.75
$ws.Range("L12").Value = 4
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("Q12").Value = 2.05
$ws.Range("R12").Value = 1.75
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.91
$ws.Range("W12").Value = 7.5
$ws.Range("AC12").Value = 9
$ws.Range("AG12").Value = 251
$ws.Range("AP12").Value = 21
$ws.Range("BA12").Value = 81
